$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new row of data (row 19): label text in B19, hyperlink in C19
$ws.Range("B19").Style = "Normal"
$ws.Range("B19").Value = "5.1. MongoDB Atlas"

$ws.Range("C19").Style = "Hyperlink"
$ws.Hyperlinks.Add(
    $ws.Range("C19"),
    "https://github.com/nguyentienminh07102004/product-management/commit/558d759eeb8d271edd366b6f9c1191d0dca0397e"
)

# Update the active selection to match the author's final cursor position
$ws.Range("B25").Select()

$wb.Save()
